# Automated update of the PEBCOM tracking sheet.
# A new claim (Caso 6236 - San Jose 1157) was reported with a report date
# that falls BEFORE the existing row 55 (Caso 6269, 6/26/2025), so it is
# inserted as a new row 55, pushing the former rows 55-56 down to 56-57.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the row immediately above the insertion point first, then insert
# that copy at row 55. This seeds the new row with the same (default, no
# special style) formatting used by every other data row, instead of Excel
# inheriting the formatting of the row that gets pushed down.
$ws.Rows.Item(54).Copy()
$ws.Rows.Item(55).Insert()
$excel.CutCopyMode = $false

# Columns A, B, D, E and I hold values that are visually numbers/dates
# ("6236", "6/24/2025", "1", "807763095", "1") but must be stored as plain
# text, matching every other row in this sheet (openpyxl originally wrote
# these as inlineStr). Force text formatting before assigning so Excel
# doesn't silently convert them to a number / date serial.
$ws.Range("A55").NumberFormat = "@"
$ws.Range("B55").NumberFormat = "@"
$ws.Range("D55").NumberFormat = "@"
$ws.Range("E55").NumberFormat = "@"
$ws.Range("I55").NumberFormat = "@"

$ws.Range("A55").Value2 = "6236"
$ws.Range("B55").Value2 = "6/24/2025"
$ws.Range("C55").Value2 = "San Jose 1157"
$ws.Range("D55").Value2 = "1"
$ws.Range("E55").Value2 = "807763095"
$ws.Range("F55").Value2 = "PEBCOM"
$ws.Range("G55").Value2 = "Pendiente"
$ws.Range("H55").Value2 = "Picada"
$ws.Range("I55").Value2 = "1"
$ws.Range("J55").Value2 = "Cambio"
$ws.Range("K55").Value2 = "Sin equipos"
$ws.Range("L55").Value2 = "Terminal"
$ws.Range("M55").Value2 = -58.385887
$ws.Range("N55").Value2 = -34.621845

# Drop the "@" text format again so the new row ends up with the same
# (default / no explicit style) look as every other data row once the
# text-safe values have been written.
$ws.Range("A55:L55").ClearFormats()
